$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "27.501.46"
    "E2" = "  +2.12%  "
    "D3" = "1.873.59"
    "D4" = "1.015"
    "E4" = "  +0.68%  "
    "D5" = "313.53"
    "E5" = "  +1.35%  "
    "E6" = "  +0.69%  "
    "D7" = "0.4800"
    "E7" = "  +1.63%  "
    "D8" = "0.3772"
    "E8" = "  +2.56%  "
    "D9" = "0.07387"
    "E9" = "  +2.79%  "
    "D10" = "0.9411"
    "E10" = "  +2.00%  "
    "D11" = "20.72"
    "E11" = "  +5.73%  "
    "D12" = "0.07871"
    "E12" = "  +3.53%  "
    "D13" = "1.863.51"
    "E13" = "  +1.23%  "
    "D14" = "5.453"
    "E14" = "  +2.73%  "
    "D15" = "6.610"
    "D16" = "91.12"
    "E16" = "  +3.05%  "
    "D17" = "1.015"
    "E17" = "  +0.58%  "
    "D18" = "0.000008994"
    "E18" = "  +4.09%  "
    "E19" = "  +0.69%  "
    "D20" = "14.96"
    "E20" = "  +2.82%  "
    "D21" = "27.535.92"
    "E21" = "  +2.18%  "
    "D22" = "5.145"
    "E22" = "  +2.40%  "
    "D23" = "10.72"
    "E23" = "  +0.68%  "
    "D24" = "1.950"
    "E24" = "  +1.78%  "
    "D25" = "153.99"
    "E25" = "  +1.14%  "
    "D26" = "18.58"
    "E26" = "  +2.31%  "
    "D27" = "2.024"
    "E27" = "  +0.78%  "
    "D28" = "116.21"
    "E28" = "  +1.63%  "
    "D29" = "5.006"
    "E29" = "  +2.40%  "
    "D30" = "0.08927"
    "E30" = "  +0.81%  "
    "D31" = "3.328"
    "E31" = "  +1.51%  "
    "D32" = "1.220"
    "E32" = "  +4.78%  "
    "D33" = "4.608"
    "E33" = "  +2.61%  "
    "D34" = "0.7519"
    "E34" = "  +0.50%  "
    "D35" = "2.712"
    "E35" = "  -2.48%  "
    "E36" = "  +6.47%  "
    "D37" = "1.123"
    "E37" = "  +2.98%  "
    "D38" = "0.05313"
    "E38" = "  +1.00%  "
    "D39" = "3.004"
    "E39" = "  +1.07%  "
    "D40" = "0.5387"
    "E40" = "  +3.45%  "
    "D41" = "7.104"
    "E41" = "  +2.93%  "
    "D42" = "0.1525"
    "E42" = "  +0.90%  "
    "D43" = "8.451"
    "E43" = "  +3.03%  "
    "D44" = "0.4852"
    "E44" = "  +3.31%  "
    "D45" = "10.64"
    "E45" = "  +1.62%  "
    "E46" = "  +0.77%  "
    "D47" = "1.665"
    "E47" = "  +3.76%  "
    "D48" = "103.19"
    "E48" = "  +1.12%  "
    "D49" = "67.16"
    "E49" = "  +2.75%  "
    "D50" = "0.06110"
    "E50" = "  +1.43%  "
    "D51" = "0.9016"
    "E51" = "  +1.85%  "
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
    $cell.Style = "Normal"
}
